$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting the existing rows 92-136 down
# to 93-137 (formatting/styles carry over from the row below, matching
# the D-column date style already used throughout the table).
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly price record.
$ws.Cells.Item(92, 1).Value = 10
$ws.Cells.Item(92, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(92, 3).Value = "La Araucanía"
$ws.Cells.Item(92, 4).Value = "2022-09-21"
$ws.Cells.Item(92, 5).Value = 9
$ws.Cells.Item(92, 6).Value = 100112035
$ws.Cells.Item(92, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 40
$ws.Cells.Item(92, 11).Value = 24000
$ws.Cells.Item(92, 12).Value = 24000
$ws.Cells.Item(92, 13).Value = 24000
$ws.Cells.Item(92, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 2400
$ws.Cells.Item(92, 17).Value = 10
$ws.Cells.Item(92, 18).Value = "Hortaliza"
